$wb = $excel.ActiveWorkbook

# Add a new worksheet. Excel.Application inserts new sheets immediately
# before the active sheet, so this lands before "ODI Batting" (the only
# existing sheet, and the active one), making it the first sheet overall.
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

# Populate the header row with bold styling (matches the "ODI Batting" header style).
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
$playerInfo.Range("A1:D1").Font.Bold = $true
$playerInfo.Range("A1:D1").HorizontalAlignment = -4108
$playerInfo.Range("A1:D1").VerticalAlignment = -4160
$playerInfo.Range("A1:D1").Borders.LineStyle = 1

# Populate the data row. The ID column looks numeric, so force it to be
# stored as text (matching the rest of the workbook, where every cell,
# including numeric-looking ones, is stored as text) before assigning it.
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "3808"
$playerInfo.Range("B2").Value = "Oliver James Hairs"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# Re-fetch the "ODI Batting" sheet by name (its index shifted after the
# insert above) and update it: rename MATCH_CARD_LINK -> MATCH_CODE and
# replace the full URLs with just the numeric match code (kept as text).
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("D1").Value = "MATCH_CODE"
$odiBatting.Range("D2:D6").NumberFormat = "@"
$odiBatting.Range("D2").Value = "3137"
$odiBatting.Range("D3").Value = "3139"
$odiBatting.Range("D4").Value = "3143"
$odiBatting.Range("D5").Value = "3146"
$odiBatting.Range("D6").Value = "3149"
